$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking "Price" values such as 0.999 / 7.21 / 0.0000240 read as plain
# numbers to Excel's type inference, which would silently normalise them
# (drop trailing zeros, switch to scientific notation, etc.) and break the
# exact text the site renders. Mark just those cells as Text first so the
# literal characters survive, then write every updated cell below.
$textGuardCells = @("D4", "D5", "D6", "D8", "D11", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D40", "D42", "D43", "D45", "D46", "D48", "D49", "D51")
foreach ($addr in $textGuardCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.125.42"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "3.102.20"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.64%  "
$ws.Range("D5").Value = "593.09"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "156.97"
$ws.Range("E6").Value = "  +3.22%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "0.540"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "3.100.55"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("D11").Value = "5.93"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  -2.61%  "
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  -3.35%  "
$ws.Range("D14").Value = "37.03"
$ws.Range("E14").Value = "  -4.07%  "
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "3.611.46"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "7.21"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "63.968.09"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "3.099.75"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").Value = "480.87"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "14.50"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").Value = "0.714"
$ws.Range("D23").Value = "7.58"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("D24").Value = "2.49"
$ws.Range("E24").Value = "  +3.85%  "
$ws.Range("D25").Value = "81.56"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "12.96"
$ws.Range("E26").Value = "  -3.97%  "
$ws.Range("D27").Value = "10.79"
$ws.Range("E27").Value = "  +9.04%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "7.62"
$ws.Range("E29").Value = "  +3.48%  "
$ws.Range("D30").Value = "2.69"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").Value = "2.20"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("D34").Value = "27.27"
$ws.Range("D35").Value = "0.0₃0845"
$ws.Range("E35").Value = "  -3.68%  "
$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").Value = "6.04"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("D38").Value = "2.27"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("E39").Value = "  -5.74%  "
$ws.Range("D40").Value = "51.10"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").Value = "443.87"
$ws.Range("E42").Value = "  -4.23%  "
$ws.Range("D43").Value = "0.292"
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("D45").Value = "0.0365"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("D46").Value = "40.38"
$ws.Range("E46").Value = "  +6.36%  "
$ws.Range("D47").Value = "2.837.95"
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("D48").Value = "131.79"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "26.28"
$ws.Range("E49").Value = "  +2.52%  "
$ws.Range("D51").Value = "2.25"
$ws.Range("E51").Value = "  -1.74%  "
